{"js": "// Remove \"Struts\" (and its trailing \", \") from the tech-stack lists.\n// e.g. \"... CSS, Struts, SVN.\" -> \"... CSS, SVN.\"\n// and   \"... d3.js, Struts, Tomcat, ...\" -> \"... d3.js, Tomcat, ...\"\nconst body = context.document.body;\n\n// Search for the exact substring \"Struts, \" (word + comma + space) so the\n// surrounding punctuation collapses back to a single, correctly formatted\n// separator list once the match is deleted.\nconst results = body.search(\"Struts, \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Remove \"Struts\" (and its trailing \", \") from the tech-stack lists.\n# e.g. \"... CSS, Struts, SVN.\" -> \"... CSS, SVN.\"\n# and  \"... d3.js, Struts, Tomcat, ...\" -> \"... d3.js, Tomcat, ...\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Struts, \"\n$find.Replacement.Text = \"\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue: keep searching the whole story\n\n# wdReplaceAll = 2 -> replace every match of \"Struts, \" with \"\" in one pass.\n$find.Execute(\n    \"Struts, \",   # FindText\n    $true,        # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"\",           # ReplaceWith\n    2             # Replace (wdReplaceAll)\n) | Out-Null\n"}
